# data/regions.xlsx -- "Add files via upload"
#
# The three counties that used to share the label "szász szék" for their
# C column (Beszterce vidék / Brassó vidék / Naszód vidék -- all part of
# Erdély) are relabeled to the new term "szász vidék". Setting a fresh
# string value on these cells causes Excel to append a new shared-string
# entry rather than reuse the old one (the old "szász szék" string is
# still used for other rows, e.g. row 18 "Fogaras vidék").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = "szász vidék"
$ws.Range("C17").Value = "szász vidék"
$ws.Range("C58").Value = "szász vidék"

# Widen column A so the (now longer) county names fit.
$ws.Columns.Item(1).ColumnWidth = 17.5

# Scroll the view back to the top and move the active selection to C12.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("C12").Select()
